$d = $word.ActiveDocument

function Find-ParagraphIndex($marker) {
    $idx = 0
    foreach ($p in $d.Paragraphs) {
        $idx = $idx + 1
        if ($p.Range.Text.StartsWith($marker)) {
            return $idx
        }
    }
    return -1
}

function Set-CleanParagraphText($marker, $newText) {
    # Remove the whole paragraph (content + paragraph mark), which merges it
    # into the following paragraph, then insert a brand-new paragraph break
    # before that point and fill it with a single clean run of text. This
    # gets rid of the mid-sentence run splits / proofErr (grammar-check)
    # markers that were littering the original runs.
    $paraIndex = Find-ParagraphIndex $marker
    $p = $d.Paragraphs($paraIndex)
    $full = $p.Range
    $full.Delete()
    $np = $d.Paragraphs($paraIndex)
    $np.Range.InsertParagraphBefore()
    $target = $d.Paragraphs($paraIndex).Range
    $target.Text = $newText
}

Set-CleanParagraphText "REQ-1.9.1" "REQ-1.9.1    The web scraper application will use the following colors for the primary elements of the user interface in Light mode: #212529, #6c757d, #adb5bd, #f8f9fa, #ffffff"
Set-CleanParagraphText "REQ-1.9.2" "REQ-1.9.2    The web scraper application will use the following colors for the primary elements of the user interface in Dark mode: #121212, #1e1e1e, #2c2c2c, #f8f9fa, #ced4da, #495057, #66b2ff"
Set-CleanParagraphText "REQ-1.9.3" "REQ-1.9.3    The web scraper application will use the following colors for the primary elements of the user interface in Blue mode: #e3f2fd, #bbdefb, #90caf9, #0d47a1, #1976d2, #1565c0, #0d47a1"
Set-CleanParagraphText "REQ-1.9.4" "REQ-1.9.4    The web scraper application will use the following colors for the primary elements of the user interface in Disco mode: #ff00ff, #00ced1, #c8bca7, #4b0082, #8d6e63, #6d4c41, #7fff00, #5d40037"

# The three blank paragraphs right after REQ-1.9.4 (and before the closing
# "Requirements that will likely need to change..." note) used to be empty;
# fill in the middle one with the new general requirement called out in the
# commit message ("Added one last general requirement").
$disco4 = Find-ParagraphIndex "REQ-1.9.4"
$d.Paragraphs($disco4 + 2).Range.Text = "REQ-1.10.1    Any user data handled by the web scraper must be minimized, securely stored, and disposed of once it is no longer needed."
